$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The sheet holds weekly price records for "Acelga" (chard) sold at the
# Terminal La Palmera de La Serena market, two rows per week (Primera /
# Segunda quality), ordered from most-recent week at the top. A new week
# of data is being added at the top of the existing run (rows 384-385),
# which pushes every following record down by two rows.
# ---------------------------------------------------------------------------

# Insert two blank rows right after the current row 385 (i.e. before the
# current row 386). Excel automatically inherits the formatting of the
# row above (including the date number format on column D).
$ws.Rows.Item(386).Resize(2).EntireRow.Insert()

# The rows that used to be 384 and 385 (now pushed down to 386 and 387)
# lost their values during the insert, so copy them back from 384/385.
$ws.Range("A384:R384").Copy($ws.Range("A386:R386"))
$ws.Range("A385:R385").Copy($ws.Range("A387:R387"))

# Now overwrite rows 384/385 with the new week's data.

# Row 384 - Primera
$ws.Cells.Item(384, 1).Value = 8
$ws.Cells.Item(384, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(384, 3).Value = "Coquimbo"
$ws.Cells.Item(384, 4).Value = 44706
$ws.Cells.Item(384, 5).Value = 4
$ws.Cells.Item(384, 6).Value = 100112009
$ws.Cells.Item(384, 7).Value = "Acelga"
$ws.Cells.Item(384, 8).Value = "Sin especificar"
$ws.Cells.Item(384, 9).Value = "Primera"
$ws.Cells.Item(384, 10).Value = 2540
$ws.Cells.Item(384, 11).Value = 600
$ws.Cells.Item(384, 12).Value = 700
$ws.Cells.Item(384, 13).Value = 650
$ws.Cells.Item(384, 14).Value = "`$/atado 1,5 a 2 kilos"
$ws.Cells.Item(384, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(384, 16).Value = 325
$ws.Cells.Item(384, 17).Value = 2
$ws.Cells.Item(384, 18).Value = "Hortaliza"

# Row 385 - Segunda
$ws.Cells.Item(385, 1).Value = 8
$ws.Cells.Item(385, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(385, 3).Value = "Coquimbo"
$ws.Cells.Item(385, 4).Value = 44706
$ws.Cells.Item(385, 5).Value = 4
$ws.Cells.Item(385, 6).Value = 100112009
$ws.Cells.Item(385, 7).Value = "Acelga"
$ws.Cells.Item(385, 8).Value = "Sin especificar"
$ws.Cells.Item(385, 9).Value = "Segunda"
$ws.Cells.Item(385, 10).Value = 1360
$ws.Cells.Item(385, 11).Value = 200
$ws.Cells.Item(385, 12).Value = 550
$ws.Cells.Item(385, 13).Value = 375
$ws.Cells.Item(385, 14).Value = "`$/atado 1,5 a 2 kilos"
$ws.Cells.Item(385, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(385, 16).Value = 188
$ws.Cells.Item(385, 17).Value = 2
$ws.Cells.Item(385, 18).Value = "Hortaliza"

Write-Output "done"
